# Switching to winter DST time
# Shift all timestamps in column A (rows 2-97) forward by 7 days
# (the data now covers the following week), and refresh the
# "Actual Production (MW)" values in column B (rows 2-97) with the
# newly fetched data for that week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 97

# Shift each existing timestamp forward by exactly 7 days, preserving
# the time-of-day fraction.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2()
    $cell.Value = $serial + 7
}

# New wind production values (MW) for rows 2..97.
$bValues = @(949,953,950,944,911,914,911,910,892,890,897,941,970,1025,1097,1178,1244,1300,1327,1358,1393,1389,1422,1514,1556,1591,1608,1644,1733,1768,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 2).Value = $bValues[$i]
}
